$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the value 12 from A2 to A3
$ws.Range("A2").ClearContents()
$ws.Range("A3").Value = 12

# Update the active selection shown in the sheet view
$ws.Range("D14").Select()
